$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.495.00'
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("D3").Value = '3.561.97'
$ws.Range("E3").Value = '  +7.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '637.48'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.49'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +9.75%  '
$ws.Range("E8").Value = '  +5.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.03'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +12.18%  '
$ws.Range("D11").Value = '3.556.36'
$ws.Range("E11").Value = '  +7.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.96'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +6.44%  '
$ws.Range("E13").Value = '  +5.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.32'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +6.81%  '
$ws.Range("D15").Value = '4.232.77'
$ws.Range("E15").Value = '  +7.51%  '
$ws.Range("D16").Value = '95.291.47'
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000254'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.43'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +6.24%  '
$ws.Range("D19").Value = '3.562.76'
$ws.Range("E19").Value = '  +7.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.97'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +19.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.17'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +6.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.508'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +13.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '517.36'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.86%  '
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.78'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +13.14%  '
$ws.Range("E26").Value = '  +9.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '93.11'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.17%  '
$ws.Range("E28").Value = '  +7.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.06'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +16.82%  '
$ws.Range("E30").Value = '  +7.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.67'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +6.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.184'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.992'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '30.34'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.570'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '588.15'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +13.31%  '
$ws.Range("E38").Value = '  +5.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +8.20%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.152'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.932'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +8.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.74'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("E44").Value = '  +6.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.85'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.57'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E48").Value = '  +5.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.13'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.23'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.14'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.05%  '
